$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: MonsterIds -> EnemyIds
$ws.Range("F1").Value = "EnemyIds"

# Row 2 (ChapterId 1, WaveId 1)
$ws.Range("F2").Value = "1,2,3,4"
$ws.Range("G2").Value = "1,3"
$ws.Range("H2").Value = "2,6"

# Row 3 (ChapterId 1, WaveId 2)
$ws.Range("F3").Value = "2,3,4,5"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "1,2,3"

# Row 4 (ChapterId 1, WaveId 3)
$ws.Range("F4").Value = "1,2,3"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "1,2"

# Rows 5-34 all share the same EnemyIds/FrontIds/BackIds pattern
for ($r = 5; $r -le 34; $r++) {
    $ws.Range("F$r").Value = "1,1,1,1"
    $ws.Range("G$r").Value = 1
    $ws.Range("H$r").Value = "1,2,3"
}

# Match the saved selection state
$null = $ws.Range("H3").Select()
